$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("21CRB00868", "Bunner", "CRIMINAL DAMAGING PROPERTY", "2909.06(A)", "M2", "Dismissed"),
    @("21TRD09437", "Bunner", "DUS", "4510.11", "M1", "Not Guilty"),
    @("21TRD09437", "Bunner", "1ST SPEED 1 YR SCHOOL >35MPHM4", "4511.21B1A", "M4", "Not Guilty"),
    @("21TRD09437", "Bunner", "RECKLESS OPERATION 1ST IN 1 YR", "4511.20", "MM", "Not Guilty"),
    @("21CRB01268", "Bunner", "POSSESSION DRUG PARAPHERNALIA", "2925.14(C)", "M4", "Not Guilty"),
    @("21CRB01268", "Bunner", "POSSESSION DRUG PARAPHERNALIA", "2925.14(C)", "M4", "Not Guilty"),
    @("21CRB01268", "Bunner", "POSSESSION DRUG PARAPHERNALIA", "2925.14(C)", "M4", "Not Guilty")
)

$startRow = 112
$endRow = $startRow + $data.Length - 1
$newRange = $ws.Range("A$startRow`:F$endRow")

# Temporarily force the new block to Text so that numeric-looking values
# like "4510.11" / "4511.20" are stored as strings rather than being
# coerced into floating point numbers. Re-applying the Normal style
# afterwards drops the temporary formatting again so the cells end up
# with the same (default) style as the rest of the sheet.
$newRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}

$newRange.Style = "Normal"
